$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.034.62"
$ws.Range("E2").Value = "  +4.71%  "
$ws.Range("D3").Value = "2.233.74"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "251.59"
$ws.Range("E5").Value = "  +7.04%  "
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").Value = "75.32"
$ws.Range("E7").Value = "  +9.09%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +5.33%  "
$ws.Range("D10").Value = "41.12"
$ws.Range("E10").Value = "  +6.86%  "
$ws.Range("D11").Value = "0.0923"
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("D12").Value = "6.87"
$ws.Range("E12").Value = "  +4.86%  "
$ws.Range("E13").Value = "  +2.70%  "
$ws.Range("D14").Value = "2.571.45"
$ws.Range("E14").Value = "  +4.80%  "
$ws.Range("D15").Value = "14.51"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "2.234.06"
$ws.Range("E16").Value = "  +4.34%  "
$ws.Range("D17").Value = "0.789"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").Value = "42.918.89"
$ws.Range("E18").Value = "  +4.92%  "
$ws.Range("E19").Value = "  +4.64%  "
$ws.Range("D20").Value = "71.17"
$ws.Range("D21").Value = "5.95"
$ws.Range("E21").Value = "  +4.98%  "
$ws.Range("D22").Value = "229.65"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.36%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "10.68"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "3.41"
$ws.Range("E27").Value = "  +4.66%  "
$ws.Range("D28").Value = "39.42"
$ws.Range("E28").Value = "  +28.86%  "
$ws.Range("E29").Value = "  +3.89%  "
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "20.18"
$ws.Range("E32").Value = "  +3.56%  "
$ws.Range("D33").Value = "0.0794"
$ws.Range("E33").Value = "  +6.86%  "
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("E35").Value = "  +9.96%  "
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").Value = "4.49"
$ws.Range("E37").Value = "  +8.84%  "
$ws.Range("D38").Value = "0.0327"
$ws.Range("E38").Value = "  +14.17%  "
$ws.Range("D39").Value = "12.32"
$ws.Range("E39").Value = "  +6.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.63%  "
$ws.Range("E41").Value = "  +10.55%  "
$ws.Range("D42").Value = "5.37"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("D43").Value = "59.56"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("B44").Value = "WOONetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D44").Value = "0.488"
$ws.Range("E44").Value = "  +32.66%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "8.64"
$ws.Range("E45").Value = "  +5.88%  "
$ws.Range("D46").Value = "103.13"
$ws.Range("E46").Value = "  +7.03%  "
$ws.Range("D47").Value = "0.0985"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("E48").Value = "  +14.61%  "
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("E50").Value = "  +3.84%  "
$ws.Range("E51").Value = "  +2.22%  "
